# Generate Report for Handoff
# Updates the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps
# for the c544b854-5d28-4257-a468-800385d347f3.md file (row 7 on every
# sheet) to reflect the freshly generated handoff xliff files.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-28 22:41:47"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-28 22:41:43"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-28 22:41:47"
